$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 8.587800000000001
$ws.Range("B12").Value = 6.103800000000001
$ws.Range("C13").Value = -12.5997
$ws.Range("B18").Value = 4.920800000000003
